$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 59.3
$ws.Range("K3").Value = 55.3
$ws.Range("K4").Value = 53.3
$ws.Range("K5").Value = 52.3

$ws.Range("N2").Value = 51.15965480231979
$ws.Range("N3").Value = 51.15965480231979
$ws.Range("N4").Value = 51.15965480231979
$ws.Range("N5").Value = 51.15965480231979
